$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the workbook window size (bookViews/workbookView)
$excel.Width = 30240
$excel.Height = 11500

# Update all A2:A104 values from 7 to 10
for ($r = 2; $r -le 104; $r++) {
    $ws.Cells.Item($r, 1).Value = 10
}

# Scroll so that row 103 is the top-left visible row, and set the new selection
$ws.Application.ActiveWindow.ScrollRow = 103
$ws.Range("A3:A104").Select()
